$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
    @('D2', '265.09'),
    @('D4', '6.229'),
    @('D5', '0.06166'),
    @('D6', '3.561'),
    @('D8', '1.360'),
    @('D9', '0.8139'),
    @('B10', 'WazirX'),
    @('C10', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @('D10', '0.1606'),
    @('E10', '9WazirXWRX'),
    @('B11', 'MandalaExchangeToken'),
    @('C11', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @('D11', '0.08214'),
    @('E11', '10MandalaExchangeTokenMDX'),
    @('B12', 'LiechtensteinCryptoassetsExchange'),
    @('C12', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
    @('D12', '0.03376'),
    @('E12', '11LiechtensteinCryptoassetsExchangeLCX'),
    @('B13', 'BitrueCoin'),
    @('C13', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @('D13', '0.03145'),
    @('E13', '12BitrueCoinBTR'),
    @('B14', 'BitMartToken'),
    @('C14', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @('D14', '0.09231'),
    @('E14', '13BitMartTokenBMX'),
    @('B15', 'MCDex'),
    @('C15', 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'),
    @('D15', '3.911'),
    @('E15', '14MCDexMCB'),
    @('B16', 'BitForexToken'),
    @('C16', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @('D16', '0.001702'),
    @('E16', '15BitForexTokenBF'),
    @('B17', 'CoinExToken'),
    @('C17', 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'),
    @('D17', '0.04839'),
    @('E17', '16CoinExTokenCET'),
    @('B18', 'One'),
    @('C18', 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'),
    @('D18', '0.0006258'),
    @('E18', '17OneONEWorstin24h'),
    @('D19', '0.006207'),
    @('D20', '0.006272'),
    @('D21', '0.001097'),
    @('D23', '3.701'),
    @('D26', '0.1246'),
    @('D27', '0.0002680'),
    @('D41', '0.007273'),
    @('D42', '0.1134'),
    @('D43', '0.003227'),
    @('D44', '0.01042'),
    @('D45', '0.00006158'),
    @('D47', '0.7697'),
    @('D48', '0.2012'),
    @('D50', '0.01240'),
)

foreach ($pair in $cells) {
    $ref = $pair[0]
    $val = $pair[1]
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}
